$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.743.70'

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.108.83'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +5.44%  '

$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.08%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '333.16'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +3.13%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.000'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.05%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5293'

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4354'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +5.02%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.08959'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +2.81%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '47.35'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +10.77%  '

$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +2.89%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '24.73'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -0.32%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '2.106.30'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +5.35%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.723'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +2.92%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.763'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +4.30%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '96.77'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +2.82%  '

$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +0.05%  '

$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +1.28%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06687'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +2.10%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '19.02'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +0.65%  '

$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +0.09%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.303'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +2.81%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '30.810.78'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +1.61%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '12.26'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +4.10%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.351.83'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +5.33%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.281'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +2.64%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '22.58'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +0.18%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.576'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +6.98%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '162.25'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -0.93%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '132.88'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +1.08%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.192'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +4.70%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.1080'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +2.76%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.167'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +1.52%  '

$ws.Range('B34').NumberFormat = '@'
$ws.Range('B34').Value = 'ARBITRUM'
$ws.Range('C34').NumberFormat = '@'
$ws.Range('C34').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.548'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +15.86%  '

$ws.Range('B35').NumberFormat = '@'
$ws.Range('B35').Value = 'HuobiToken'
$ws.Range('C35').NumberFormat = '@'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.896'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +1.67%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02591'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +3.15%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '9.585'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +6.63%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.535'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +3.01%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.06756'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +2.59%  '

$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +3.73%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.2270'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +3.13%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.6835'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +3.20%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.245'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +1.39%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.9999'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -0.12%  '

$ws.Range('B45').NumberFormat = '@'
$ws.Range('B45').Value = 'Decentraland'
$ws.Range('C45').NumberFormat = '@'
$ws.Range('C45').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.6414'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +4.19%  '

$ws.Range('B46').NumberFormat = '@'
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').NumberFormat = '@'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '14.03'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +2.34%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.222'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +1.06%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.652'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -0.11%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.260'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -0.74%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '82.98'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +3.80%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.191'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +7.62%  '
